# Error Calculations and Plots
# Remove the "RM 232" and "SC 92" rows entirely (the whole row shifts up),
# then update the simulated-missing-data pattern in the remaining rows to
# match the new random mask (some previously blanked cells get their
# numeric values restored, other cells get newly blanked).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely - all following rows shift up by one.
$ws.Rows(26).Delete()

# After the above deletion, the row that used to be "SC 92" (originally
# row 28) is now row 27. Delete it too.
$ws.Rows(27).Delete()

# Restore / clear individual data cells to match the new missing pattern.
$ws.Range("D2").Value = -13.5
$ws.Range("E3").ClearContents()
$ws.Range("E4").Value = -6.4
$ws.Range("D6").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("D12").Value = -14.1
$ws.Range("D14").ClearContents()
$ws.Range("E15").Value = -8.4
$ws.Range("E18").Value = -8.5
$ws.Range("E19").ClearContents()
$ws.Range("D20").Value = -14
$ws.Range("D21").Value = -14.3
$ws.Range("E22").ClearContents()
$ws.Range("D23").ClearContents()
$ws.Range("E23").Value = -7
$ws.Range("D24").ClearContents()
$ws.Range("E25").Value = -7.1
$ws.Range("B26").Value = -20.2
$ws.Range("B27").ClearContents()
$ws.Range("E27").ClearContents()
$ws.Range("B28").ClearContents()
$ws.Range("B29").Value = -19.5
$ws.Range("B30").Value = -19.7
$ws.Range("B31").ClearContents()
$ws.Range("D31").Value = -13.7
$ws.Range("B32").ClearContents()
$ws.Range("D33").Value = -14.1
